$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.999.30'
$ws.Range('E2').Value = '  -3.30%  '
$ws.Range('D3').Value = '1.681.39'
$ws.Range('E3').Value = '  -3.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.59'
$ws.Range('E5').Value = '  -1.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9983'
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3670'
$ws.Range('E7').Value = '  -2.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3374'
$ws.Range('E8').Value = '  -6.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.57'
$ws.Range('E9').Value = '  -5.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.176'
$ws.Range('E10').Value = '  -4.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07331'
$ws.Range('E11').Value = '  -3.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9973'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.191'
$ws.Range('E13').Value = '  -3.19%  '
$ws.Range('E14').Value = '  -5.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.828'
$ws.Range('E15').Value = '  -3.11%  '
$ws.Range('D16').Value = '1.676.98'
$ws.Range('E16').Value = '  -3.34%  '
$ws.Range('E17').Value = '  -4.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06605'
$ws.Range('E18').Value = '  -2.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9977'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '82.25'
$ws.Range('E20').Value = '  -4.51%  '
$ws.Range('E21').Value = '  -3.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.194'
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.64'
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('D24').Value = '24.871.76'
$ws.Range('E24').Value = '  -3.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.431'
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.698'
$ws.Range('E26').Value = '  -4.99%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.84'
$ws.Range('E27').Value = '  -2.92%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '150.45'
$ws.Range('E28').Value = '  -3.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.268'
$ws.Range('E29').Value = '  +8.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '130.20'
$ws.Range('E30').Value = '  -2.95%  '
$ws.Range('D31').Value = '1.864.22'
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.507'
$ws.Range('E32').Value = '  -5.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.153'
$ws.Range('E33').Value = '  +1.17%  '
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08600'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.733'
$ws.Range('E36').Value = '  -3.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.439'
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06482'
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02346'
$ws.Range('E39').Value = '  -5.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.736'
$ws.Range('E40').Value = '  -5.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2168'
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.248'
$ws.Range('E42').Value = '  -3.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6274'
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9973'
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.43'
$ws.Range('E45').Value = '  -2.22%  '
$ws.Range('E46').Value = '  -2.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5986'
$ws.Range('E47').Value = '  -4.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.045'
$ws.Range('E48').Value = '  -4.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.88'
$ws.Range('E49').Value = '  -4.68%  '
$ws.Range('E50').Value = '  -3.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '77.37'
$ws.Range('E51').Value = '  -2.00%  '
